$d = $word.ActiveDocument

# Locate the paragraph whose text is "CONNECT<nbsp>:" (the nbsp is U+00A0).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like ("CONNECT" + [char]0xA0 + ":*")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'CONNECT :' paragraph"
}

$r = $target.Range

# Rebuild the paragraph: pStyle -> Heading3, and strip the direct
# character/paragraph-mark formatting (the fr-FR/zxx/zxx w:lang override)
# so both w:pPr/w:rPr and w:r/w:rPr come back empty, matching a plain
# "Heading 3" styled line.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="Heading3"/><w:rPr/></w:pPr>' +
       '<w:r><w:rPr/><w:t>CONNECT' + [char]0xA0 + ':</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

[void]$r.InsertXML($xml)
